# Applies cell-value updates to the cryptocurrency price/volume table
# described by the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.430.19"
$ws.Range("E2").Value = "  -2.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.276.46"
$ws.Range("E3").Value = "  -4.52%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.40"
$ws.Range("E5").Value = "  -3.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.05"
$ws.Range("E6").Value = "  -7.05%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -4.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.78"
$ws.Range("E10").Value = "  -6.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "50.68"
$ws.Range("E11").Value = "  -5.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0787"
$ws.Range("E12").Value = "  -3.20%  "
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.65"
$ws.Range("E14").Value = "  -4.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.626.62"
$ws.Range("E15").Value = "  -4.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.18"
$ws.Range("E16").Value = "  -2.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.261.78"
$ws.Range("E17").Value = "  -5.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.785"
$ws.Range("E18").Value = "  -3.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.342.70"
$ws.Range("E19").Value = "  -2.85%  "
$ws.Range("E20").Value = "  -2.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.44"
$ws.Range("E21").Value = "  -3.92%  "
$ws.Range("E22").Value = "  -5.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.59"
$ws.Range("E23").Value = "  -2.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "234.58"
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("E25").Value = "  -5.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.49"
$ws.Range("E26").Value = "  -5.06%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.34"
$ws.Range("E28").Value = "  -6.09%  "
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.85"
$ws.Range("E30").Value = "  -7.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "163.87"
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.08"
$ws.Range("E32").Value = "  -4.68%  "
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.95"
$ws.Range("E34").Value = "  -5.70%  "
$ws.Range("E35").Value = "  -4.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0697"
$ws.Range("E36").Value = "  -5.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.36"
$ws.Range("E37").Value = "  -6.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.83"
$ws.Range("E38").Value = "  -9.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.10"
$ws.Range("E39").Value = "  -12.46%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.76"
$ws.Range("E40").Value = "  -9.01%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0998"
$ws.Range("E41").Value = "  -5.83%  "
$ws.Range("E42").Value = "  -3.41%  "
$ws.Range("E43").Value = "  -7.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.965.91"
$ws.Range("E44").Value = "  -3.51%  "
$ws.Range("E45").Value = "  -3.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.81"
$ws.Range("E46").Value = "  -9.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.64"
$ws.Range("E47").Value = "  -8.73%  "
$ws.Range("E48").Value = "  -9.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.82"
$ws.Range("E49").Value = "  -4.63%  "
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.500.84"
$ws.Range("E51").Value = "  -4.16%  "
